# The commit swaps the two theme parts of the deck: the design that is
# actually applied to the slide master / slides (currently the "Integral"
# theme, physically stored as ppt/theme/theme2.xml) is replaced by the
# stock "Office Theme" palette (previously stored as ppt/theme/theme1.xml),
# while the other (unused, notes-master-only) theme part swaps the other
# way. The PowerPoint object model only exposes the live, applied theme
# through Design/SlideMaster.Theme.ThemeColorScheme, so we drive the
# visible change - the 12 theme colours that cascade through every
# schemeClr reference on every slide - through that API.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

# Target palette: the stock Office theme colour scheme (was theme1.xml,
# becomes the colour scheme behind the presentation's single design).
# ThemeColorScheme.Colors(n) is 1-based and follows clrScheme document
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink. RGB is written
# as a standard BGR-packed OLE color (R + G*256 + B*65536).
$themeColors.Colors(1).RGB = 0          # dk1      000000
$themeColors.Colors(2).RGB = 16777215   # lt1      FFFFFF
$themeColors.Colors(3).RGB = 6968388    # dk2      44546A
$themeColors.Colors(4).RGB = 15132391   # lt2      E7E6E6
$themeColors.Colors(5).RGB = 13998939   # accent1  5B9BD5
$themeColors.Colors(6).RGB = 3243501    # accent2  ED7D31
$themeColors.Colors(7).RGB = 10855845   # accent3  A5A5A5
$themeColors.Colors(8).RGB = 49407      # accent4  FFC000
$themeColors.Colors(9).RGB = 12874308   # accent5  4472C4
$themeColors.Colors(10).RGB = 4697456   # accent6  70AD47
$themeColors.Colors(11).RGB = 12673797  # hlink    0563C1
$themeColors.Colors(12).RGB = 7491477   # folHlink 954F72
